$d = $word.ActiveDocument

# 1. Remove trailing space from "Expected graduation: Spring 2027 "
$d.Content.Find.Execute(
    "Expected graduation: Spring 2027 ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Expected graduation: Spring 2027", 2) | Out-Null

# 2. Remove the appended " CGPA: 3.7/4.0" text from the specialization line,
#    keeping the existing trailing double space.
$d.Content.Find.Execute(
    "Specialization: Econometrics, Energy & Resource Economics, Industrial Organization  CGPA: 3.7/4.0 ",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Specialization: Econometrics, Energy & Resource Economics, Industrial Organization  ", 2) | Out-Null

# 3. Delete the standalone "CGPA: 3.8/4.0" paragraph entirely (text + paragraph mark).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "CGPA: 3.8/4.0 `r") {
        $p.Range.Delete()
        break
    }
}

# 4. Delete the standalone "Achieved first class" paragraph entirely (text + paragraph mark).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Achieved first class `r") {
        $p.Range.Delete()
        break
    }
}

# 5. Remove only the run text "CGPA: 3.3/4.0 " but keep the (now empty) paragraph.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "CGPA: 3.3/4.0 `r") {
        $r = $p.Range
        $r.MoveEnd(1, -1) | Out-Null
        $r.Text = ""
        break
    }
}
